$wb = $excel.ActiveWorkbook

# --- Sheet 1: "People & Projects" ---
# Row 3 ("institution"/"funding" demo values) and Row 4 ("middle2"/"email2"/"project2")
# lost their sample data -- clear those cells.
$wsPeople = $wb.Worksheets.Item("People & Projects")
$wsPeople.Range("F3").ClearContents()
$wsPeople.Range("G3").ClearContents()
$wsPeople.Range("B4").ClearContents()
$wsPeople.Range("E4").ClearContents()
$wsPeople.Range("H4").ClearContents()

# --- Sheet 3: "Data Types" ---
# Row 3 ("samplemethod") and Row 4 ("datamethod2") lost their sample data.
$wsDataTypes = $wb.Worksheets.Item("Data Types")
$wsDataTypes.Range("F3").ClearContents()
$wsDataTypes.Range("G4").ClearContents()

# --- Selections / active sheet ---
# Final selection on "People & Projects" moves to E4.
[void]$wsPeople.Range("E4").Select()

# "Data Types" becomes the active/visible tab (was "Package Descriptions").
[void]$wsDataTypes.Activate()
